$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 566, shifting rows 566:589 down to 568:591.
$ws.Rows("566:567").Insert()

# Populate new row 566 with the new weekly price entry (Morada(o), Primera)
$ws.Range("A566").Value = 11
$ws.Range("B566").Value = "Vega Monumental Concepción"
$ws.Range("C566").Value = "Bíobío"
$ws.Range("D566").Value = 45106
$ws.Range("E566").Value = 8
$ws.Range("F566").Value = 100112006
$ws.Range("G566").Value = "Repollo"
$ws.Range("H566").Value = "Morada(o)"
$ws.Range("I566").Value = "Primera"
$ws.Range("J566").Value = 1000
$ws.Range("K566").Value = 800
$ws.Range("L566").Value = 900
$ws.Range("M566").Value = 850
$ws.Range("N566").Value = "$/unidad"
$ws.Range("O566").Value = "Región Metropolitana"
$ws.Range("P566").Value = 850
$ws.Range("Q566").Value = 1
$ws.Range("R566").Value = "Hortaliza"

# Populate new row 567 with the new weekly price entry (Morada(o), Segunda)
$ws.Range("A567").Value = 11
$ws.Range("B567").Value = "Vega Monumental Concepción"
$ws.Range("C567").Value = "Bíobío"
$ws.Range("D567").Value = 45106
$ws.Range("E567").Value = 8
$ws.Range("F567").Value = 100112006
$ws.Range("G567").Value = "Repollo"
$ws.Range("H567").Value = "Morada(o)"
$ws.Range("I567").Value = "Segunda"
$ws.Range("J567").Value = 500
$ws.Range("K567").Value = 700
$ws.Range("L567").Value = 700
$ws.Range("M567").Value = 700
$ws.Range("N567").Value = "$/unidad"
$ws.Range("O567").Value = "Región Metropolitana"
$ws.Range("P567").Value = 700
$ws.Range("Q567").Value = 1
$ws.Range("R567").Value = "Hortaliza"
